$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.479.92"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.572.29"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.49"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3753"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3412"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07595"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.43"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.969"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "1.569.90"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.25"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.285"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.47"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.24"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").Value = "22.462.13"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.327"
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.613"
$ws.Range("E26").Value = "  -5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.72"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.996"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.30"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "1.746.71"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.214"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.933"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08463"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.379"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02476"
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2302"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06566"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.504"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.47"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6307"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.10"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.819"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5904"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.27"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.230"
$ws.Range("E50").Value = "  -5.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07334"
$ws.Range("E51").Value = "  -0.12%  "
